$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of "Random" method results appended below the existing data (rows 2-19).
$rows = @(
    @(42602.58185185185,  "Random", 0, 0, 0, 0, 0, 65, 35, 0, 0, 95, 5),
    @(42602.97996527778,  "Random", 0, 0, 0, 0, 0, 46, 54, 0, 0, 8,  92),
    @(42603.694537037038, "Random", 0, 0, 0, 0, 0, 53, 47, 0, 0, 1,  99)
)

$startRow = 20
$lastExistingRow = $startRow - 1

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    # Copy the date cell's existing style (from the prior row) so no new
    # number format / style entry gets created - just reuse style index 1.
    $ws.Cells.Item($lastExistingRow, 1).Copy() | Out-Null
    $ws.Cells.Item($r, 1).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Cells.Item($r, 1).Value = $data[0]

    $ws.Cells.Item($r, 2).Value = $data[1]

    for ($c = 3; $c -le 13; $c++) {
        $ws.Cells.Item($r, $c).Value = $data[$c - 1]
    }
}
